$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly-data block (row 500), pushing the
# existing rows 500:604 down to 502:606 (the last two rows that fall off the
# bottom reappear as the new rows 605:606 automatically).
$ws.Rows("500:501").Insert()

# New row 500 - "Primera" quality, week of 2022-03-21 (serial 44641)
$ws.Range("A500").Value = 6
$ws.Range("B500").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C500").Value = "Metropolitana"
$ws.Range("D500").Value = 44641
$ws.Range("E500").Value = 13
$ws.Range("F500").Value = 100112017
$ws.Range("G500").Value = "Apio"
$ws.Range("H500").Value = "Americana (o)"
$ws.Range("I500").Value = "Primera"
$ws.Range("J500").Value = 2010
$ws.Range("K500").Value = 7000
$ws.Range("L500").Value = 8000
$ws.Range("M500").Value = 7517
$ws.Range("N500").Value = "`$/docena de matas"
$ws.Range("O500").Value = "Región de Coquimbo"
$ws.Range("P500").Value = 1253
$ws.Range("Q500").Value = 6
$ws.Range("R500").Value = "Hortaliza"

# New row 501 - "Segunda" quality, same week
$ws.Range("A501").Value = 6
$ws.Range("B501").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C501").Value = "Metropolitana"
$ws.Range("D501").Value = 44641
$ws.Range("E501").Value = 13
$ws.Range("F501").Value = 100112017
$ws.Range("G501").Value = "Apio"
$ws.Range("H501").Value = "Americana (o)"
$ws.Range("I501").Value = "Segunda"
$ws.Range("J501").Value = 620
$ws.Range("K501").Value = 5000
$ws.Range("L501").Value = 6000
$ws.Range("M501").Value = 5677
$ws.Range("N501").Value = "`$/docena de matas"
$ws.Range("O501").Value = "Región de Coquimbo"
$ws.Range("P501").Value = 946
$ws.Range("Q501").Value = 6
$ws.Range("R501").Value = "Hortaliza"
